$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.502.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.78%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.470.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9591"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "281.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.96%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3703"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3177"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.057"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06676"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("E12").Value = "  -0.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.604"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.80%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.251"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.47%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.474.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001036"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05721"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9579"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.71%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.669"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.256"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.680.26"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.291"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "137.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.637.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.956"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.319"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8323"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +29.37%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07835"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.27%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.915"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02073"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9703"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.27%  "

# Row 41
$ws.Range("E41").Value = "  +2.56%  "

# Row 42
$ws.Range("E42").Value = "  -1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.324"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5408"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.56%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.591"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.65%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.42%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5316"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.57%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.818"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "

# Row 50
$ws.Range("E50").Value = "  +4.43%  "

# Row 51
$ws.Range("E51").Value = "  -0.06%  "
